# Rename the worksheets to reflect the cleaned-up dataset structure:
#   "upstream"          -> "emissions"
#   "upstream inflows"  -> "removals"
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("upstream").Name = "emissions"
$wb.Worksheets.Item("upstream inflows").Name = "removals"
